# Update RMSE, MAE, PAEM, Bias, IA, CE, R2 validation metrics
# for Mean & conditional mean imputation results (rows 2-25, columns E-K)

$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("E2").Value = 0.177040614553816
$ws.Range("F2").Value = 0.03271345881050904
$ws.Range("G2").Value = 0.4119776363568582
$ws.Range("H2").Value = 0.008006391738267952
$ws.Range("I2").Value = 0.9842737936058645
$ws.Range("J2").Value = 0.9399032630736774
$ws.Range("K2").Value = 0.9400261706758157
$ws.Range("E3").Value = 0.2504929342755191
$ws.Range("F3").Value = 0.06579994006005384
$ws.Range("G3").Value = 0.9075282800135882
$ws.Range("H3").Value = 0.01763694502351926
$ws.Range("I3").Value = 0.9670114716162004
$ws.Range("J3").Value = 0.8796915767461381
$ws.Range("K3").Value = 0.8802879961950569
$ws.Range("E4").Value = 0.350764879548896
$ws.Range("F4").Value = 0.1287472743760802
$ws.Range("G4").Value = 1.953999776834536
$ws.Range("H4").Value = 0.03797411871229358
$ws.Range("I4").Value = 0.9290822933606959
$ws.Range("J4").Value = 0.7640949266995968
$ws.Range("K4").Value = 0.7668598334082359
$ws.Range("E5").Value = 0.5026378515062118
$ws.Range("F5").Value = 0.2578747337779352
$ws.Range("G5").Value = 4.555505770754011
$ws.Range("H5").Value = 0.08853190209335432
$ws.Range("I5").Value = 0.8213773232463539
$ws.Range("J5").Value = 0.5155873726728482
$ws.Range("K5").Value = 0.5306154930423113
$ws.Range("E6").Value = 0.1660203616698156
$ws.Range("F6").Value = 0.0310521315129072
$ws.Range("G6").Value = 0.05602688982520403
$ws.Range("H6").Value = -0.001088829072820848
$ws.Range("I6").Value = 0.9862432539284421
$ws.Range("J6").Value = 0.9471520937351563
$ws.Range("K6").Value = 0.9471543668663025
$ws.Range("E7").Value = 0.2274589130094877
$ws.Range("F7").Value = 0.05977118312009332
$ws.Range("G7").Value = 0.05155052184960606
$ws.Range("H7").Value = -0.001001835138164088
$ws.Range("I7").Value = 0.9732083559213293
$ws.Range("J7").Value = 0.9008001549960689
$ws.Range("K7").Value = 0.9008020794060287
$ws.Range("E8").Value = 0.3173328993539681
$ws.Range("F8").Value = 0.1176754147957396
$ws.Range("G8").Value = 0.07103903785795875
$ws.Range("H8").Value = -0.001380575826469851
$ws.Range("I8").Value = 0.9435673203049786
$ws.Range("J8").Value = 0.8069208962233845
$ws.Range("K8").Value = 0.8069245507035201
$ws.Range("E9").Value = 0.4572649771465216
$ws.Range("F9").Value = 0.2398180566095743
$ws.Range("G9").Value = 0.2584969347059904
$ws.Range("H9").Value = 0.005023640944929605
$ws.Range("I9").Value = 0.8569303371359109
$ws.Range("J9").Value = 0.5990954796412684
$ws.Range("K9").Value = 0.5991438681522158
$ws.Range("E10").Value = 0.14841338187132
$ws.Range("F10").Value = 0.02896461124791324
$ws.Range("G10").Value = 1.128162121046377
$ws.Range("H10").Value = -0.02192475291923015
$ws.Range("I10").Value = 0.989118351435694
$ws.Range("J10").Value = 0.9577670714574064
$ws.Range("K10").Value = 0.9586887394174771
$ws.Range("E11").Value = 0.2090171111843771
$ws.Range("F11").Value = 0.05791753302368297
$ws.Range("G11").Value = 2.262589300265958
$ws.Range("H11").Value = -0.0439712612581022
$ws.Range("I11").Value = 0.9778366781061348
$ws.Range("J11").Value = 0.9162338110769621
$ws.Range("K11").Value = 0.9199409846516621
$ws.Range("E12").Value = 0.3113067625788649
$ws.Range("F12").Value = 0.1214975840605807
$ws.Range("G12").Value = 4.831866850497923
$ws.Range("H12").Value = -0.09390271562878576
$ws.Range("I12").Value = 0.9480491465028942
$ws.Range("J12").Value = 0.8141843944919857
$ws.Range("K12").Value = 0.8310911930815574
$ws.Range("E13").Value = 0.4932480576010425
$ws.Range("F13").Value = 0.2756586299876386
$ws.Range("G13").Value = 11.77618746608658
$ws.Range("H13").Value = -0.2288589518366477
$ws.Range("I13").Value = 0.8592237844598907
$ws.Range("J13").Value = 0.5335169776962454
$ws.Range("K13").Value = 0.63394175652331
$ws.Range("E14").Value = 0.1282826116788999
$ws.Range("F14").Value = 0.02331397681915248
$ws.Range("G14").Value = 0.5255336684444168
$ws.Range("H14").Value = -0.01571896525981276
$ws.Range("I14").Value = 0.9862191523888855
$ws.Range("J14").Value = 0.9469721082430758
$ws.Range("K14").Value = 0.9477682981782146
$ws.Range("E15").Value = 0.1854237224963548
$ws.Range("F15").Value = 0.04713375850113642
$ws.Range("G15").Value = 1.001296417951039
$ws.Range("H15").Value = -0.02994925835129821
$ws.Range("I15").Value = 0.9700023675051184
$ws.Range("J15").Value = 0.8892103051690939
$ws.Range("K15").Value = 0.8921005918356066
$ws.Range("E16").Value = 0.2626065651168993
$ws.Range("F16").Value = 0.09350780578339697
$ws.Range("G16").Value = 2.127143742047527
$ws.Range("H16").Value = -0.0636238943222143
$ws.Range("I16").Value = 0.9354540316014135
$ws.Range("J16").Value = 0.777781642446197
$ws.Range("K16").Value = 0.7908255968245903
$ws.Range("E17").Value = 0.3951294098703838
$ws.Range("F17").Value = 0.199361938009991
$ws.Range("G17").Value = 5.051942597345183
$ws.Range("H17").Value = -0.1511060374443661
$ws.Range("I17").Value = 0.8347385769282054
$ws.Range("J17").Value = 0.4969079127510794
$ws.Range("K17").Value = 0.5704832749700395
$ws.Range("E18").Value = 0.1221380329802444
$ws.Range("F18").Value = 0.02187045800782594
$ws.Range("G18").Value = 0.01882188585003358
$ws.Range("H18").Value = -0.0005629716753953907
$ws.Range("I18").Value = 0.9875334151433026
$ws.Range("J18").Value = 0.9519303878654223
$ws.Range("K18").Value = 0.9519314091384667
$ws.Range("E19").Value = 0.1755159288102985
$ws.Range("F19").Value = 0.0450456767318077
$ws.Range("G19").Value = 0.06897052873495001
$ws.Range("H19").Value = -0.002062941748993312
$ws.Range("I19").Value = 0.973192132977132
$ws.Range("J19").Value = 0.9007336968531102
$ws.Range("K19").Value = 0.9007474101761161
$ws.Range("E20").Value = 0.2487180774885371
$ws.Range("F20").Value = 0.08989904407627372
$ws.Range("G20").Value = 0.04203717288874267
$ws.Range("H20").Value = 0.001257352097373332
$ws.Range("I20").Value = 0.9414194850983194
$ws.Range("J20").Value = 0.8006650373173626
$ws.Range("K20").Value = 0.8006701315948495
$ws.Range("E21").Value = 0.3493049891969582
$ws.Range("F21").Value = 0.1783333631414632
$ws.Range("G21").Value = 0.05672445803967791
$ws.Range("H21").Value = -0.0016966558735365
$ws.Range("I21").Value = 0.8607018598245459
$ws.Range("J21").Value = 0.6068318286357915
$ws.Range("K21").Value = 0.6068411045410677
$ws.Range("E22").Value = 0.1294871590151127
$ws.Range("F22").Value = 0.02333603148897975
$ws.Range("G22").Value = 0.5594626062284366
$ws.Range("H22").Value = -0.01673379613812357
$ws.Range("I22").Value = 0.9859510106246359
$ws.Range("J22").Value = 0.9459715908755906
$ws.Range("K22").Value = 0.9468739049314069
$ws.Range("E23").Value = 0.1825755512972048
$ws.Range("F23").Value = 0.04649145484630868
$ws.Range("G23").Value = 1.105726015480153
$ws.Range("H23").Value = -0.03307279793443338
$ws.Range("I23").Value = 0.9710307711694485
$ws.Range("J23").Value = 0.8925876999439615
$ws.Range("K23").Value = 0.8961123065422846
$ws.Range("E24").Value = 0.2663383968565355
$ws.Range("F24").Value = 0.09543493023529215
$ws.Range("G24").Value = 2.439108068645621
$ws.Range("H24").Value = -0.07295489765566668
$ws.Range("I24").Value = 0.9337840851572148
$ws.Range("J24").Value = 0.7714209951520811
$ws.Range("K24").Value = 0.7885715304664422
$ws.Range("E25").Value = 0.409150627249522
$ws.Range("F25").Value = 0.2076560953838874
$ws.Range("G25").Value = 5.813976309088061
$ws.Range("H25").Value = -0.1738988329604909
$ws.Range("I25").Value = 0.8250723763962009
$ws.Range("J25").Value = 0.4605698486428372
$ws.Range("K25").Value = 0.5580154207913334
